$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AX97").Value = 'Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning'

$ws.Range("B98").Value = 98935
$ws.Range("AX98").Value = 'Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning'

$ws.Range("B99").Value = 98935
$ws.Range("AX99").Value = 'Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning'

$ws.Range("AX100").Value = 'Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning'

$ws.Range("A102").Value = 130964547
$ws.Range("B102").Value = 57881
$ws.Range("E102").Value = 100049
$ws.Range("F102").Value = 'Spillkråka'
$ws.Range("G102").Value = 'Dryocopus martius'
$ws.Range("H102").Value = '(Linnaeus, 1758)'
$ws.Range("Q102").Value = 509495
$ws.Range("R102").Value = 6718877
$ws.Range("AC102").Value = 'Födosökspår . inventering åt vasa vind'

$ws.Range("A103").Value = 130964526
$ws.Range("B103").Value = 79245
$ws.Range("E103").Value = 6425
$ws.Range("F103").Value = 'Garnlav'
$ws.Range("G103").Value = 'Alectoria sarmentosa'
$ws.Range("H103").Value = '(Ach.) Ach.'
$ws.Range("Q103").Value = 509610
$ws.Range("R103").Value = 6719050
$ws.Range("AC103").Value = 'Enstaka . inventering åt vasa vind'
$ws.Range("AX103").Value = 'Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning'

$ws.Range("A104").Value = 130964545
$ws.Range("B104").Value = 57073
$ws.Range("D104").Value = 'LC'
$ws.Range("E104").Value = 100138
$ws.Range("F104").Value = 'Tjäder'
$ws.Range("G104").Value = 'Tetrao urogallus'
$ws.Range("H104").Value = 'Linnaeus, 1758'
$ws.Range("Q104").Value = 509535
$ws.Range("R104").Value = 6718925
$ws.Range("AC104").Value = 'Spillning . inventering åt vasa vind'
$ws.Range("AX104").Value = 'Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning'

$ws.Range("A105").Value = 130964541
$ws.Range("B105").Value = 91813
$ws.Range("D105").Value = 'NT'
$ws.Range("E105").Value = 1202
$ws.Range("F105").Value = 'Ullticka'
$ws.Range("G105").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H105").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q105").Value = 509703
$ws.Range("R105").Value = 6719018
$ws.Range("AC105").Value = 'Enstaka . inventering åt vasa vind'
$ws.Range("AX105").Value = 'Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning'

$ws.Range("A106").Value = 130964537
$ws.Range("B106").Value = 79245
$ws.Range("D106").Value = 'NT'
$ws.Range("E106").Value = 6425
$ws.Range("F106").Value = 'Garnlav'
$ws.Range("G106").Value = 'Alectoria sarmentosa'
$ws.Range("H106").Value = '(Ach.) Ach.'
$ws.Range("Q106").Value = 509822
$ws.Range("R106").Value = 6718960
$ws.Range("AC106").Value = 'Rikligt . inventering åt vasa vind'

$ws.Range("A107").Value = 130964642
$ws.Range("B107").Value = 99041
$ws.Range("D107").Value = 'LC'
$ws.Range("E107").Value = 221952
$ws.Range("F107").Value = 'Spindelblomster'
$ws.Range("G107").Value = 'Neottia cordata'
$ws.Range("H107").Value = '(L.) Rich.'
$ws.Range("Q107").Value = 509917
$ws.Range("R107").Value = 6719042
$ws.Range("AC107").Value = 'Måttliga förekomster . inventering åt vasa vind'
$ws.Range("AX107").Value = 'Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning'

$ws.Range("AX108").Value = 'Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning'

$ws.Range("AX109").Value = 'Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning'

$ws.Range("B110").Value = 98935

$ws.Range("B111").Value = 92111

$ws.Range("A112").Value = 130964650
$ws.Range("B112").Value = 92272
$ws.Range("D112").Value = 'VU'
$ws.Range("E112").Value = 1209
$ws.Range("F112").Value = 'Rynkskinn'
$ws.Range("G112").Value = 'Hermanssonia centrifuga'
$ws.Range("H112").Value = '(P. Karst.) Zmitr.'
$ws.Range("Q112").Value = 509694
$ws.Range("R112").Value = 6718936

$ws.Range("A114").Value = 130964645
$ws.Range("B114").Value = 99041
$ws.Range("D114").Value = 'LC'
$ws.Range("E114").Value = 221952
$ws.Range("F114").Value = 'Spindelblomster'
$ws.Range("G114").Value = 'Neottia cordata'
$ws.Range("H114").Value = '(L.) Rich.'
$ws.Range("Q114").Value = 509804
$ws.Range("R114").Value = 6719024

$ws.Range("B116").Value = 98935
$ws.Range("AX116").Value = 'Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning'

$ws.Range("B118").Value = 98935

$ws.Range("B119").Value = 92272

$ws.Range("B120").Value = 98935

$ws.Range("AX121").Value = 'Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning'

$ws.Range("A122").Value = 130964640
$ws.Range("B122").Value = 57881
$ws.Range("D122").Value = 'NT'
$ws.Range("E122").Value = 100049
$ws.Range("F122").Value = 'Spillkråka'
$ws.Range("G122").Value = 'Dryocopus martius'
$ws.Range("H122").Value = '(Linnaeus, 1758)'
$ws.Range("Q122").Value = 509697
$ws.Range("R122").Value = 6719144
$ws.Range("AC122").Value = 'Gamla födosöksspår . inventering åt vasa vind'
$ws.Range("AX122").Value = 'Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning'

$ws.Range("A123").Value = 130964542
$ws.Range("B123").Value = 57073
$ws.Range("D123").Value = 'LC'
$ws.Range("E123").Value = 100138
$ws.Range("F123").Value = 'Tjäder'
$ws.Range("G123").Value = 'Tetrao urogallus'
$ws.Range("H123").Value = 'Linnaeus, 1758'
$ws.Range("Q123").Value = 509635
$ws.Range("R123").Value = 6718941
$ws.Range("AC123").Value = 'Spillning . inventering åt vasa vind'

$ws.Range("B124").Value = 98922

$ws.Range("B125").Value = 98935

$ws.Range("B126").Value = 99018
$ws.Range("AX126").Value = 'Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning'

$ws.Range("A127").Value = 130964546
$ws.Range("B127").Value = 92508
$ws.Range("D127").Value = 'VU'
$ws.Range("E127").Value = 898
$ws.Range("F127").Value = 'Blackticka'
$ws.Range("G127").Value = 'Steccherinum collabens'
$ws.Range("H127").Value = '(Fr.) Vesterholt'
$ws.Range("Q127").Value = 509515
$ws.Range("R127").Value = 6718886
$ws.Range("AC127").Value = 'Betydande förekomst . inventering åt vasa vind'
$ws.Range("AX127").Value = 'Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning'

$ws.Range("A128").Value = 130964538
$ws.Range("B128").Value = 79245
$ws.Range("D128").Value = 'NT'
$ws.Range("E128").Value = 6425
$ws.Range("F128").Value = 'Garnlav'
$ws.Range("G128").Value = 'Alectoria sarmentosa'
$ws.Range("H128").Value = '(Ach.) Ach.'
$ws.Range("Q128").Value = 509875
$ws.Range("R128").Value = 6719025
$ws.Range("AC128").Value = 'Enstaka . inventering åt vasa vind'
$ws.Range("AX128").Value = 'Sofia Berg, Pia Edfors, Anna Sjövall, Anders Esplund, Enviro Planning'

$ws.Range("B129").Value = 91813
